$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows appended to the bottom of the data (rows 6-8)
$data = @(
    @(9994.8799999999992, 9975.93, 78.05, 78.2, $false, 0.19, 42613.766319444447, $true),
    @(10076.84, 9994.8799999999992, 77.739999999999995, 78.38, $false, 0.82, 42614.674351851849, $true),
    @(10022.43, 10076.84, 78.36, 77.94, $false, -0.54, 42615.752500000002, $false)
)

$rowIndex = 6
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}

# Column G holds dates - copy the existing date formatting (style index) down
# onto the new rows so the new cells share the same number format as G3:G5.
$ws.Range("G5").Copy()
$ws.Range("G6:G8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The extra rows push the bottom of column A's best-fit width out slightly.
$ws.Columns.Item(1).ColumnWidth = 8.14
